# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.624.03"
$ws.Range("E2").Value = "  +5.68%  "

$ws.Range("D3").Value = "3.184.37"
$ws.Range("E3").Value = "  +2.88%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "401.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.77"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.25%  "

$ws.Range("E11").Value = "  +1.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0886"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.94%  "

$ws.Range("D13").Value = "3.675.66"
$ws.Range("E13").Value = "  +2.70%  "

$ws.Range("E14").Value = "  +1.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.06%  "

$ws.Range("E16").Value = "  +8.81%  "

$ws.Range("D17").Value = "3.185.46"
$ws.Range("E17").Value = "  +3.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.52"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("D19").Value = "54.500.26"
$ws.Range("E19").Value = "  +5.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.06%  "

$ws.Range("E21").Value = "  +3.35%  "

$ws.Range("D22").Value = "0.0₃0997"
$ws.Range("E22").Value = "  +3.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.28"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.06"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.45"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.42%  "

$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("E31").Value = "  +4.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.06"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.61%  "

$ws.Range("E33").Value = "  +13.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "37.06"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.09"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.88"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("E37").Value = "  +7.46%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("E39").Value = "  +10.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.17"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +12.71%  "

$ws.Range("E41").Value = "  +3.05%  "

$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.68%  "

$ws.Range("E45").Value = "  +1.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("D49").Value = "2.091.52"
$ws.Range("E49").Value = "  +2.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0348"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +9.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0508"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +11.52%  "
